$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68
$ws.Range("A68").Value = 67.0
$ws.Range("B68").Value = 'Tuesday, Jan 10'
$ws.Range("C68").Value = '6:09 AM'
$ws.Range("D68").Value = 'P81988'
$ws.Range("E68").Value = 'Cologne'
$ws.Range("F68").Value = '(CGN)'
$ws.Range("G68").Value = 'SprintAir '
$ws.Range("H68").Value = 'AT72'
$ws.Range("I68").Value = '(SP-SPG)'
$ws.Range("J68").Value = '6:14 AM'
$ws.Range("L68").Value = '0 hours, 5 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(68, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(68, 13))

# Row 69
$ws.Range("A69").Value = 68.0
$ws.Range("B69").Value = 'Tuesday, Jan 10'
$ws.Range("C69").Value = '8:05 AM'
$ws.Range("D69").Value = 'FR4999'
$ws.Range("E69").Value = 'Brussels'
$ws.Range("F69").Value = '(CRL)'
$ws.Range("G69").Value = 'Ryanair '
$ws.Range("H69").Value = 'B738'
$ws.Range("I69").Value = '(EI-DPG)'
$ws.Range("J69").Value = '8:09 AM'
$ws.Range("L69").Value = '0 hours, 4 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(69, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(69, 13))

# Row 70
$ws.Range("A70").Value = 69.0
$ws.Range("B70").Value = 'Tuesday, Jan 10'
$ws.Range("C70").Value = '10:15 AM'
$ws.Range("D70").Value = 'FR7897'
$ws.Range("E70").Value = 'Paris'
$ws.Range("F70").Value = '(BVA)'
$ws.Range("G70").Value = 'Ryanair '
$ws.Range("H70").Value = 'B738'
$ws.Range("I70").Value = '(SP-RSM)'
$ws.Range("J70").Value = '10:02 AM'
$ws.Range("L70").Value = '0 hours, -13 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(70, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(70, 13))

# Row 71
$ws.Range("A71").Value = 70.0
$ws.Range("B71").Value = 'Tuesday, Jan 10'
$ws.Range("C71").Value = '11:20 AM'
$ws.Range("D71").Value = 'LO3945'
$ws.Range("E71").Value = 'Warsaw'
$ws.Range("F71").Value = '(WAW)'
$ws.Range("G71").Value = 'LOT '
$ws.Range("H71").Value = 'E75S'
$ws.Range("I71").Value = '(SP-LIB)'
$ws.Range("J71").Value = '11:30 AM'
$ws.Range("L71").Value = '0 hours, 10 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(71, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(71, 13))

# Row 72
$ws.Range("A72").Value = 71.0
$ws.Range("B72").Value = 'Tuesday, Jan 10'
$ws.Range("C72").Value = '11:20 AM'
$ws.Range("D72").Value = 'LO6595'
$ws.Range("E72").Value = 'Warsaw'
$ws.Range("F72").Value = '(WAW)'
$ws.Range("G72").Value = 'LOT '
$ws.Range("H72").Value = 'B788'
$ws.Range("I72").Value = '(SP-LRD)'
$ws.Range("J72").Value = '11:21 AM'
$ws.Range("L72").Value = '0 hours, 1 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(72, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(72, 13))

# Row 73
$ws.Range("A73").Value = 72.0
$ws.Range("B73").Value = 'Tuesday, Jan 10'
$ws.Range("C73").Value = '12:20 PM'
$ws.Range("D73").Value = 'LH1636'
$ws.Range("E73").Value = 'Munich'
$ws.Range("F73").Value = '(MUC)'
$ws.Range("G73").Value = 'Lufthansa '
$ws.Range("H73").Value = 'CRJ9'
$ws.Range("I73").Value = '(D-ACNN)'
$ws.Range("J73").Value = '12:12 PM'
$ws.Range("L73").Value = '0 hours, -8 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(73, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(73, 13))

# Row 74
$ws.Range("A74").Value = 73.0
$ws.Range("B74").Value = 'Tuesday, Jan 10'
$ws.Range("C74").Value = '1:45 PM'
$ws.Range("D74").Value = 'FR3727'
$ws.Range("E74").Value = 'Billund'
$ws.Range("F74").Value = '(BLL)'
$ws.Range("G74").Value = 'Ryanair '
$ws.Range("H74").Value = 'B738'
$ws.Range("I74").Value = '(SP-RSM)'
$ws.Range("J74").Value = '1:30 PM'
$ws.Range("L74").Value = '0 hours, -15 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(74, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(74, 13))

# Row 75
$ws.Range("A75").Value = 74.0
$ws.Range("B75").Value = 'Tuesday, Jan 10'
$ws.Range("C75").Value = '2:15 PM'
$ws.Range("D75").Value = 'LO3943'
$ws.Range("E75").Value = 'Warsaw'
$ws.Range("F75").Value = '(WAW)'
$ws.Range("G75").Value = 'LOT '
$ws.Range("H75").Value = 'E75S'
$ws.Range("I75").Value = '(SP-LIA)'
$ws.Range("J75").Value = '2:10 PM'
$ws.Range("L75").Value = '0 hours, -5 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(75, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(75, 13))

# Row 76
$ws.Range("A76").Value = 75.0
$ws.Range("B76").Value = 'Tuesday, Jan 10'
$ws.Range("C76").Value = '2:50 PM'
$ws.Range("D76").Value = 'LH1390'
$ws.Range("E76").Value = 'Frankfurt'
$ws.Range("F76").Value = '(FRA)'
$ws.Range("G76").Value = 'Lufthansa '
$ws.Range("H76").Value = 'CRJ9'
$ws.Range("I76").Value = '(D-ACNJ)'
$ws.Range("J76").Value = '2:55 PM'
$ws.Range("L76").Value = '0 hours, 5 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(76, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(76, 13))

# Row 77
$ws.Range("A77").Value = 76.0
$ws.Range("B77").Value = 'Tuesday, Jan 10'
$ws.Range("C77").Value = '3:01 PM'
$ws.Range("D77").Value = 'UNKNOWN'
$ws.Range("E77").Value = 'Poprad'
$ws.Range("F77").Value = '(TAT)'
$ws.Range("G77").Value = 'AMC Aviation '
$ws.Range("H77").Value = 'PC24'
$ws.Range("I77").Value = '(SP-AGA)'
$ws.Range("J77").Value = '2:52 PM'
$ws.Range("L77").Value = '0 hours, -9 minutes'
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(77, 11))
$ws.Cells.Item(2, 13).Copy($ws.Cells.Item(77, 13))
